# Registration.Samples.xlsx — add the "Object handles" factory sample section
# (rows 22-31) to Sheet1, matching the target commit:
#   "An example of how one could implement object handles"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: the overall AND() now needs to cover the new tests too ---------
$ws.Range("B1").Formula = "=AND(A4:A31)"

# --- Section header (row 22) -------------------------------------------
# Shared-string insertion order matters (it drives the new sst indices), so
# we deliberately write these in the same order the original authoring
# session must have: D22, G22, then the "One"/"Two" labels as the table
# rows are filled in, and finally the A22 row caption.
$ws.Range("D22").Value = "One call returning IEnumerable<T>"
$ws.Range("G22").Value = "Multiple calls returning T"

# --- Table block 1: rows 23-26 -------------------------------------------
$ws.Range("D23:D25").FormulaArray = "=_xll.dnaFactoryMultiple(E23:E25,F23:F25)"

$ws.Range("E23").Value = "One"
$ws.Range("F23").Value = 1
$ws.Range("G23").Formula = "=_xll.dnaFactorySingle(E23,F23)"

$ws.Range("E24").Value = "Two"
$ws.Range("F24").Value = 2
$ws.Range("G24").Formula = "=_xll.dnaFactorySingle(E24,F24)"
$ws.Range("J24").Formula = "=_xll.dnaFactoryCompound(G24,G23)"

$ws.Range("E25").Value = "One"
$ws.Range("F25").Value = 3
$ws.Range("G25").Formula = "=_xll.dnaFactorySingle(E25,F25)"
$ws.Range("J25").Formula = "=_xll.dnaFactoryCompound(J24,G25)"

$ws.Range("A26").Formula = "=SUMSQ(B26:C26)<0.00000000000001"
$ws.Range("B26").Formula = "=D26-G26"
$ws.Range("C26").Formula = "=G26-J26"
$ws.Range("D26").Formula = "=_xll.dnaUseSomeHandles(D23:D25)"
$ws.Range("G26").Formula = "=_xll.dnaUseSomeHandles(G23:G25)"
$ws.Range("J26").Formula = "=_xll.dnaUseSomeHandles(J25)"

# --- Table block 2: rows 28-31 -------------------------------------------
$ws.Range("D28:D30").FormulaArray = "=_xll.dnaFactoryMultiple(E28:E30,F28:F30)"

$ws.Range("E28").Value = "Two"
$ws.Range("F28").Value = 4
$ws.Range("G28").Formula = "=_xll.dnaFactorySingle(E28,F28)"

$ws.Range("E29").Value = "One"
$ws.Range("F29").Value = 5
$ws.Range("G29").Formula = "=_xll.dnaFactorySingle(E29,F29)"
$ws.Range("J29").Formula = "=_xll.dnaFactoryCompound(G29,G28)"

$ws.Range("E30").Value = "Two"
$ws.Range("F30").Value = 6
$ws.Range("G30").Formula = "=_xll.dnaFactorySingle(E30,F30)"
$ws.Range("J30").Formula = "=_xll.dnaFactoryCompound(J29,G30)"

$ws.Range("A31").Formula = "=SUMSQ(B31:C31)<0.00000000000001"
$ws.Range("B31").Formula = "=D31-G31"
$ws.Range("C31").Formula = "=G31-J31"
$ws.Range("D31").Formula = "=_xll.dnaUseSomeHandles(D28:D30)"
$ws.Range("G31").Formula = "=_xll.dnaUseSomeHandles(G28:G30)"
$ws.Range("J31").Formula = "=_xll.dnaUseSomeHandles(J30)"

# Row caption, written last so its string lands at the end of the shared
# string table (matches the target sst ordering).
$ws.Range("A22").Value = "Object handles"

# --- Bold style for the three header cells --------------------------------
$ws.Range("A22").Font.Bold = $true
$ws.Range("D22").Font.Bold = $true
$ws.Range("G22").Font.Bold = $true

# --- Column D is now wide enough to show the full "One call returning
# IEnumerable<T>" caption.
$ws.Columns.Item(4).ColumnWidth = 15.35

# --- View state: scrolled down, with D32 selected -------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D32").Select()
